# Auto-generated edit script: applies the cell-value changes described by the diff
# to Sheets/Hades_Profits.xlsx (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1931.7142
$ws.Range("J6").Value = 1931.7142
$ws.Range("L6").Value = 5795.142599999999
$ws.Range("N6").Value = -6019.142599999999
$ws.Range("H19").Value = 422.44446
$ws.Range("J19").Value = 538.2727
$ws.Range("L19").Value = 538.2727
$ws.Range("N19").Value = -888.2727
$ws.Range("H43").Value = 1685.409
$ws.Range("I43").Value = 844.75
$ws.Range("K43").Value = 844.75
$ws.Range("M43").Value = -775.75
$ws.Range("H129").Value = 894.8723
$ws.Range("I129").Value = 720.82355
$ws.Range("J129").Value = 993.5
$ws.Range("K129").Value = 2162.47065
$ws.Range("L129").Value = 2980.5
$ws.Range("M129").Value = 2837.52935
$ws.Range("N129").Value = -12980.5
$ws.Range("H132").Value = 1582174.9
$ws.Range("I132").Value = 1533.68
$ws.Range("J132").Value = 8168180
$ws.Range("K132").Value = 4601.04
$ws.Range("L132").Value = 24504540
$ws.Range("M132").Value = -2071.04
$ws.Range("N132").Value = -24509600
$ws.Range("H137").Value = 1725955.9
$ws.Range("I137").Value = 3334613.8
$ws.Range("J137").Value = 2393.9644
$ws.Range("K137").Value = 10003841.4
$ws.Range("L137").Value = 7181.8932
$ws.Range("M137").Value = -10001291.4
$ws.Range("N137").Value = -12281.8932

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1753.5834
$ws.Range("I2").Value = 1666.25
$ws.Range("J2").Value = 1928.25
$ws.Range("K2").Value = 1666.25
$ws.Range("L2").Value = 1928.25
$ws.Range("M2").Value = -1553.25
$ws.Range("N2").Value = -2154.25
$ws.Range("H116").Value = 1753.5834
$ws.Range("I116").Value = 1666.25
$ws.Range("J116").Value = 1928.25
$ws.Range("K116").Value = 1666.25
$ws.Range("L116").Value = 1928.25
$ws.Range("M116").Value = 627.75
$ws.Range("N116").Value = -6516.25
$ws.Range("H132").Value = 35568.56
$ws.Range("I132").Value = 26237.25
$ws.Range("J132").Value = 55213.42
$ws.Range("K132").Value = 78711.75
$ws.Range("L132").Value = 165640.26
$ws.Range("M132").Value = -76181.75
$ws.Range("N132").Value = -170700.26

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1753.5834
$ws.Range("I3").Value = 1666.25
$ws.Range("J3").Value = 1928.25
$ws.Range("K3").Value = 1666.25
$ws.Range("L3").Value = 1928.25
$ws.Range("M3").Value = -1552.25
$ws.Range("N3").Value = -2156.25
$ws.Range("H22").Value = 286.0625
$ws.Range("I22").Value = 238.46666
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 238.46666
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -65.46665999999999
$ws.Range("N22").Value = -1346
$ws.Range("H86").Value = 17093.564
$ws.Range("I86").Value = 23001.092
$ws.Range("J86").Value = 11678.333
$ws.Range("K86").Value = 23001.092
$ws.Range("L86").Value = 11678.333
$ws.Range("M86").Value = -21878.092
$ws.Range("N86").Value = -13924.333
$ws.Range("H89").Value = 17093.564
$ws.Range("I89").Value = 23001.092
$ws.Range("J89").Value = 11678.333
$ws.Range("K89").Value = 115005.46
$ws.Range("L89").Value = 58391.665
$ws.Range("M89").Value = -109389.46
$ws.Range("N89").Value = -69623.66500000001
$ws.Range("H134").Value = 1849.4615
$ws.Range("I134").Value = 1216.0303
$ws.Range("J134").Value = 5333.3335
$ws.Range("K134").Value = 3648.0909
$ws.Range("L134").Value = 16000.0005
$ws.Range("M134").Value = -1113.0909
$ws.Range("N134").Value = -21070.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 20577.627
$ws.Range("I134").Value = 1465.2444
$ws.Range("J134").Value = 82010.28999999999
$ws.Range("K134").Value = 4395.733200000001
$ws.Range("L134").Value = 246030.87
$ws.Range("M134").Value = -1860.733200000001
$ws.Range("N134").Value = -251100.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3600
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 3600
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H80").Value = 2590.0938
$ws.Range("I80").Value = 1968
$ws.Range("J80").Value = 2872.8635
$ws.Range("K80").Value = 5904
$ws.Range("L80").Value = 8618.5905
$ws.Range("M80").Value = -4968
$ws.Range("N80").Value = -10490.5905
$ws.Range("H83").Value = 2590.0938
$ws.Range("I83").Value = 1968
$ws.Range("J83").Value = 2872.8635
$ws.Range("K83").Value = 17712
$ws.Range("L83").Value = 25855.7715
$ws.Range("M83").Value = -13032
$ws.Range("N83").Value = -35215.7715
$ws.Range("H87").Value = 12624.75
$ws.Range("I87").Value = 12624.75
$ws.Range("K87").Value = 37874.25
$ws.Range("M87").Value = -36626.25
$ws.Range("H90").Value = 12624.75
$ws.Range("I90").Value = 12624.75
$ws.Range("K90").Value = 113622.75
$ws.Range("M90").Value = -107382.75
$ws.Range("H104").Value = 3882
$ws.Range("J104").Value = 3882
$ws.Range("L104").Value = 11646
$ws.Range("N104").Value = -16888
$ws.Range("H131").Value = 14592.591
$ws.Range("J131").Value = 15965.35
$ws.Range("L131").Value = 47896.05
$ws.Range("N131").Value = -57976.05
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 6250
$ws.Range("I31").Value = 2031.25
$ws.Range("K31").Value = 2031.25
$ws.Range("M31").Value = -1739.25
$ws.Range("H37").Value = 6250
$ws.Range("I37").Value = 2031.25
$ws.Range("K37").Value = 2031.25
$ws.Range("M37").Value = -1754.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1807.9
$ws.Range("I113").Value = 1123.5
$ws.Range("J113").Value = 2264.1667
$ws.Range("K113").Value = 3370.5
$ws.Range("L113").Value = 6792.500100000001
$ws.Range("M113").Value = -1200.5
$ws.Range("N113").Value = -11132.5001
$ws.Range("H122").Value = 1804.6945
$ws.Range("I122").Value = 982.8333
$ws.Range("J122").Value = 3448.4167
$ws.Range("K122").Value = 2948.4999
$ws.Range("L122").Value = 10345.2501
$ws.Range("M122").Value = -498.4998999999998
$ws.Range("N122").Value = -15245.2501
$ws.Range("H126").Value = 1341.4667
$ws.Range("I126").Value = 1330.4762
$ws.Range("J126").Value = 1367.1111
$ws.Range("K126").Value = 3991.4286
$ws.Range("L126").Value = 4101.3333
$ws.Range("M126").Value = -1521.4286
$ws.Range("N126").Value = -9041.3333
$ws.Range("H136").Value = 34470.367
$ws.Range("I136").Value = 22245.723
$ws.Range("J136").Value = 78667.16
$ws.Range("K136").Value = 66737.169
$ws.Range("L136").Value = 236001.48
$ws.Range("M136").Value = -64187.169
$ws.Range("N136").Value = -241101.48
